$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wrap every "mecanique" label in column B (rows 2-22) in single quotes,
# e.g. "Strategie" -> "'Strategie'".
#
# A single *leading* apostrophe typed/assigned into a cell is interpreted
# by Excel (and this COM host) as a "force text" entry marker: it is
# consumed from the stored value and the cell is flagged with a
# quote-prefix style instead. To end up with a literal leading apostrophe
# in the stored text we supply an extra leading apostrophe (so one is
# consumed as the marker and one remains), then clear the resulting
# quote-prefix formatting so the cell's style stays the default (matching
# a value that was never "quote prefixed").
for ($r = 2; $r -le 22; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Text
    $cell.Value = "''" + $current + "'"
    $cell.ClearFormats()
}

# Move the active selection from A23 to F6.
$ws.Range("F6").Select()
